$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1, G1, H1 - same style as existing header (copy formatting from E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Fill F2:H21 with boolean FALSE values, matching the outlier detection flags
for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}

# The one true outlier flag: H4
$ws.Cells.Item(4, 8).Value = $true
